# Commit: Username replaced with Team
# For every sheet: rewrite the "01.xlsx"/"02.xlsx" headers to the
# full "Rank-ICPC 2025 Team Formation - N.xlsx" titles, drop the old
# "03.xlsx" column (old column D) so FinalPoints (old column E) slides
# into D, and rewrite the data rows in the new sort order / values.
$wb = $excel.ActiveWorkbook

# ---- Sheet: Participants ----
$ws = $wb.Worksheets.Item("Participants")
$ws.Range("B1").Value = "Rank-ICPC 2025 Team Formation - 01.xlsx"
$ws.Range("C1").Value = "Rank-ICPC 2025 Team Formation - 02.xlsx"
$ws.Range("D1:D17").EntireColumn.Delete()

$ws.Cells.Item(2, 1).Value = "Tamjid_Hossen(Tamjid)"
$ws.Cells.Item(2, 2).Value = 200
$ws.Cells.Item(2, 3).Value = 300
$ws.Cells.Item(2, 4).Value = 500
$ws.Cells.Item(3, 1).Value = "YouDOntKNowWHo(Nabeel Ahsan)"
$ws.Cells.Item(3, 2).Value = 225
$ws.Cells.Item(3, 3).Value = 258
$ws.Cells.Item(3, 4).Value = 483
$ws.Cells.Item(4, 1).Value = "sf61561(Syed Fahad Mahmud)"
$ws.Cells.Item(4, 2).Value = 300
$ws.Cells.Item(4, 3).Value = 164
$ws.Cells.Item(4, 4).Value = 464
$ws.Cells.Item(5, 1).Value = "shazidmashrafi(Shazid)"
$ws.Cells.Item(5, 2).Value = 180
$ws.Cells.Item(5, 3).Value = 225
$ws.Cells.Item(5, 4).Value = 405
$ws.Cells.Item(6, 1).Value = "y_this_kolaveri(SAIF)"
$ws.Cells.Item(6, 2).Value = 258
$ws.Cells.Item(6, 3).Value = 120
$ws.Cells.Item(6, 4).Value = 378
$ws.Cells.Item(7, 1).Value = "AL_AMIN_17(Al Amin)"
$ws.Cells.Item(7, 2).Value = 129
$ws.Cells.Item(7, 3).Value = 200
$ws.Cells.Item(7, 4).Value = 329
$ws.Cells.Item(8, 1).Value = "Md_Saurob_bhuyan(Noob)"
$ws.Cells.Item(8, 2).Value = 164
$ws.Cells.Item(8, 3).Value = 139
$ws.Cells.Item(8, 4).Value = 303
$ws.Cells.Item(9, 1).Value = "rakin_ahsan(Rakin)"
$ws.Cells.Item(9, 2).Value = 106
$ws.Cells.Item(9, 3).Value = 180
$ws.Cells.Item(9, 4).Value = 286
$ws.Cells.Item(10, 1).Value = "farhanshadiq(Farhan)"
$ws.Cells.Item(10, 2).Value = 113
$ws.Cells.Item(10, 3).Value = 129
$ws.Cells.Item(10, 4).Value = 242
$ws.Cells.Item(11, 1).Value = "Noornabi1770(Noor)"
$ws.Cells.Item(11, 2).Value = 90
$ws.Cells.Item(11, 3).Value = 150
$ws.Cells.Item(11, 4).Value = 240
$ws.Cells.Item(12, 1).Value = "Aniksamiul(Anik)"
$ws.Cells.Item(12, 2).Value = 120
$ws.Cells.Item(12, 3).Value = 106
$ws.Cells.Item(12, 4).Value = 226
$ws.Cells.Item(13, 1).Value = "Apon_Chowdhury(Apon)"
$ws.Cells.Item(13, 2).Value = 100
$ws.Cells.Item(13, 3).Value = 113
$ws.Cells.Item(13, 4).Value = 213
$ws.Cells.Item(14, 1).Value = "Marufhussain(maruf)"
$ws.Cells.Item(14, 2).Value = 150
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 150
$ws.Cells.Item(15, 1).Value = "Akash_khan"
$ws.Cells.Item(15, 2).Value = 139
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 139
$ws.Cells.Item(16, 1).Value = "SadmanIshtiak(Sadman)"
$ws.Cells.Item(16, 2).Value = 95
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 95
$ws.Cells.Item(17, 1).Value = "_Mohiul007(Rabby)"
$ws.Cells.Item(17, 2).Value = 86
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 86

# ---- Sheet: Team_1 ----
$ws = $wb.Worksheets.Item("Team_1")
$ws.Range("B1").Value = "Rank-ICPC 2025 Team Formation - 01.xlsx"
$ws.Range("C1").Value = "Rank-ICPC 2025 Team Formation - 02.xlsx"
$ws.Range("D1:D4").EntireColumn.Delete()

$ws.Cells.Item(2, 1).Value = "Tamjid_Hossen(Tamjid)"
$ws.Cells.Item(2, 2).Value = 200
$ws.Cells.Item(2, 3).Value = 300
$ws.Cells.Item(2, 4).Value = 500
$ws.Cells.Item(3, 1).Value = "YouDOntKNowWHo(Nabeel Ahsan)"
$ws.Cells.Item(3, 2).Value = 225
$ws.Cells.Item(3, 3).Value = 258
$ws.Cells.Item(3, 4).Value = 483
$ws.Cells.Item(4, 1).Value = "sf61561(Syed Fahad Mahmud)"
$ws.Cells.Item(4, 2).Value = 300
$ws.Cells.Item(4, 3).Value = 164
$ws.Cells.Item(4, 4).Value = 464

# ---- Sheet: Team_2 ----
$ws = $wb.Worksheets.Item("Team_2")
$ws.Range("B1").Value = "Rank-ICPC 2025 Team Formation - 01.xlsx"
$ws.Range("C1").Value = "Rank-ICPC 2025 Team Formation - 02.xlsx"
$ws.Range("D1:D4").EntireColumn.Delete()

$ws.Cells.Item(2, 1).Value = "shazidmashrafi(Shazid)"
$ws.Cells.Item(2, 2).Value = 180
$ws.Cells.Item(2, 3).Value = 225
$ws.Cells.Item(2, 4).Value = 405
$ws.Cells.Item(3, 1).Value = "y_this_kolaveri(SAIF)"
$ws.Cells.Item(3, 2).Value = 258
$ws.Cells.Item(3, 3).Value = 120
$ws.Cells.Item(3, 4).Value = 378
$ws.Cells.Item(4, 1).Value = "AL_AMIN_17(Al Amin)"
$ws.Cells.Item(4, 2).Value = 129
$ws.Cells.Item(4, 3).Value = 200
$ws.Cells.Item(4, 4).Value = 329

# ---- Sheet: Team_3 ----
$ws = $wb.Worksheets.Item("Team_3")
$ws.Range("B1").Value = "Rank-ICPC 2025 Team Formation - 01.xlsx"
$ws.Range("C1").Value = "Rank-ICPC 2025 Team Formation - 02.xlsx"
$ws.Range("D1:D4").EntireColumn.Delete()

$ws.Cells.Item(2, 1).Value = "Md_Saurob_bhuyan(Noob)"
$ws.Cells.Item(2, 2).Value = 164
$ws.Cells.Item(2, 3).Value = 139
$ws.Cells.Item(2, 4).Value = 303
$ws.Cells.Item(3, 1).Value = "rakin_ahsan(Rakin)"
$ws.Cells.Item(3, 2).Value = 106
$ws.Cells.Item(3, 3).Value = 180
$ws.Cells.Item(3, 4).Value = 286
$ws.Cells.Item(4, 1).Value = "farhanshadiq(Farhan)"
$ws.Cells.Item(4, 2).Value = 113
$ws.Cells.Item(4, 3).Value = 129
$ws.Cells.Item(4, 4).Value = 242

# ---- Sheet: Team_4 ----
$ws = $wb.Worksheets.Item("Team_4")
$ws.Range("B1").Value = "Rank-ICPC 2025 Team Formation - 01.xlsx"
$ws.Range("C1").Value = "Rank-ICPC 2025 Team Formation - 02.xlsx"
$ws.Range("D1:D4").EntireColumn.Delete()

$ws.Cells.Item(2, 1).Value = "Noornabi1770(Noor)"
$ws.Cells.Item(2, 2).Value = 90
$ws.Cells.Item(2, 3).Value = 150
$ws.Cells.Item(2, 4).Value = 240
$ws.Cells.Item(3, 1).Value = "Aniksamiul(Anik)"
$ws.Cells.Item(3, 2).Value = 120
$ws.Cells.Item(3, 3).Value = 106
$ws.Cells.Item(3, 4).Value = 226
$ws.Cells.Item(4, 1).Value = "Apon_Chowdhury(Apon)"
$ws.Cells.Item(4, 2).Value = 100
$ws.Cells.Item(4, 3).Value = 113
$ws.Cells.Item(4, 4).Value = 213

# ---- Sheet: Team_5 ----
$ws = $wb.Worksheets.Item("Team_5")
$ws.Range("B1").Value = "Rank-ICPC 2025 Team Formation - 01.xlsx"
$ws.Range("C1").Value = "Rank-ICPC 2025 Team Formation - 02.xlsx"
$ws.Range("D1:D4").EntireColumn.Delete()

$ws.Cells.Item(2, 1).Value = "Marufhussain(maruf)"
$ws.Cells.Item(2, 2).Value = 150
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 150
$ws.Cells.Item(3, 1).Value = "Akash_khan"
$ws.Cells.Item(3, 2).Value = 139
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 139
$ws.Cells.Item(4, 1).Value = "SadmanIshtiak(Sadman)"
$ws.Cells.Item(4, 2).Value = 95
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 95

# ---- Sheet: Team_6 ----
$ws = $wb.Worksheets.Item("Team_6")
$ws.Range("B1").Value = "Rank-ICPC 2025 Team Formation - 01.xlsx"
$ws.Range("C1").Value = "Rank-ICPC 2025 Team Formation - 02.xlsx"
$ws.Range("D1:D2").EntireColumn.Delete()

$ws.Cells.Item(2, 1).Value = "_Mohiul007(Rabby)"
$ws.Cells.Item(2, 2).Value = 86
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 86

Write-Host "Edit complete"